$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C77").Value = "Mean_Abs_Deviation_test"
$ws.Range("A77").Value = "Mean_Abs_Deviation"
$ws.Range("B77").Value = "Test Mean Absolute Deviation"

[void]$ws.Range("B78").Select()
